$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range('D2')
$r.NumberFormat = "@"
$r.Value = '66.634.46'
$r.Style = "Normal"
$r = $ws.Range('E2')
$r.NumberFormat = "@"
$r.Value = '  -4.35%  '
$r.Style = "Normal"
$r = $ws.Range('D3')
$r.NumberFormat = "@"
$r.Value = '3.340.06'
$r.Style = "Normal"
$r = $ws.Range('E3')
$r.NumberFormat = "@"
$r.Value = '  -1.39%  '
$r.Style = "Normal"
$r = $ws.Range('E4')
$r.NumberFormat = "@"
$r.Value = '  -0.02%  '
$r.Style = "Normal"
$r = $ws.Range('D5')
$r.NumberFormat = "@"
$r.Value = '574.00'
$r.Style = "Normal"
$r = $ws.Range('E5')
$r.NumberFormat = "@"
$r.Value = '  -3.39%  '
$r.Style = "Normal"
$r = $ws.Range('D6')
$r.NumberFormat = "@"
$r.Value = '180.65'
$r.Style = "Normal"
$r = $ws.Range('E6')
$r.NumberFormat = "@"
$r.Value = '  -5.43%  '
$r.Style = "Normal"
$r = $ws.Range('D7')
$r.NumberFormat = "@"
$r.Value = '0.628'
$r.Style = "Normal"
$r = $ws.Range('E7')
$r.NumberFormat = "@"
$r.Value = '  +3.19%  '
$r.Style = "Normal"
$r = $ws.Range('E8')
$r.NumberFormat = "@"
$r.Value = '  -0.04%  '
$r.Style = "Normal"
$r = $ws.Range('E9')
$r.NumberFormat = "@"
$r.Value = '  -3.27%  '
$r.Style = "Normal"
$r = $ws.Range('E10')
$r.NumberFormat = "@"
$r.Value = '  -1.88%  '
$r.Style = "Normal"
$r = $ws.Range('D11')
$r.NumberFormat = "@"
$r.Value = '0.404'
$r.Style = "Normal"
$r = $ws.Range('E11')
$r.NumberFormat = "@"
$r.Value = '  -3.71%  '
$r.Style = "Normal"
$r = $ws.Range('D12')
$r.NumberFormat = "@"
$r.Value = '3.916.20'
$r.Style = "Normal"
$r = $ws.Range('E12')
$r.NumberFormat = "@"
$r.Value = '  -1.58%  '
$r.Style = "Normal"
$r = $ws.Range('E13')
$r.NumberFormat = "@"
$r.Value = '  -0.30%  '
$r.Style = "Normal"
$r = $ws.Range('D14')
$r.NumberFormat = "@"
$r.Value = '27.05'
$r.Style = "Normal"
$r = $ws.Range('E14')
$r.NumberFormat = "@"
$r.Value = '  -5.89%  '
$r.Style = "Normal"
$r = $ws.Range('D15')
$r.NumberFormat = "@"
$r.Value = '66.704.09'
$r.Style = "Normal"
$r = $ws.Range('E15')
$r.NumberFormat = "@"
$r.Value = '  -4.20%  '
$r.Style = "Normal"
$r = $ws.Range('D16')
$r.NumberFormat = "@"
$r.Value = '0.0000168'
$r.Style = "Normal"
$r = $ws.Range('E16')
$r.NumberFormat = "@"
$r.Value = '  -2.51%  '
$r.Style = "Normal"
$r = $ws.Range('D17')
$r.NumberFormat = "@"
$r.Value = '3.334.62'
$r.Style = "Normal"
$r = $ws.Range('E17')
$r.NumberFormat = "@"
$r.Value = '  -2.11%  '
$r.Style = "Normal"
$r = $ws.Range('D18')
$r.NumberFormat = "@"
$r.Value = '436.50'
$r.Style = "Normal"
$r = $ws.Range('E18')
$r.NumberFormat = "@"
$r.Value = '  -3.27%  '
$r.Style = "Normal"
$r = $ws.Range('D19')
$r.NumberFormat = "@"
$r.Value = '5.69'
$r.Style = "Normal"
$r = $ws.Range('E19')
$r.NumberFormat = "@"
$r.Value = '  -2.52%  '
$r.Style = "Normal"
$r = $ws.Range('D20')
$r.NumberFormat = "@"
$r.Value = '13.61'
$r.Style = "Normal"
$r = $ws.Range('E20')
$r.NumberFormat = "@"
$r.Value = '  -1.59%  '
$r.Style = "Normal"
$r = $ws.Range('D21')
$r.NumberFormat = "@"
$r.Value = '7.61'
$r.Style = "Normal"
$r = $ws.Range('E21')
$r.NumberFormat = "@"
$r.Value = '  -2.80%  '
$r.Style = "Normal"
$r = $ws.Range('D22')
$r.NumberFormat = "@"
$r.Value = '73.55'
$r.Style = "Normal"
$r = $ws.Range('E22')
$r.NumberFormat = "@"
$r.Value = '  -3.72%  '
$r.Style = "Normal"
$r = $ws.Range('E23')
$r.NumberFormat = "@"
$r.Value = '  -0.13%  '
$r.Style = "Normal"
$r = $ws.Range('D24')
$r.NumberFormat = "@"
$r.Value = '0.518'
$r.Style = "Normal"
$r = $ws.Range('E24')
$r.NumberFormat = "@"
$r.Value = '  -0.97%  '
$r.Style = "Normal"
$r = $ws.Range('D25')
$r.NumberFormat = "@"
$r.Value = '0.0000118'
$r.Style = "Normal"
$r = $ws.Range('E25')
$r.NumberFormat = "@"
$r.Value = '  -4.10%  '
$r.Style = "Normal"
$r = $ws.Range('E26')
$r.NumberFormat = "@"
$r.Value = '  -0.12%  '
$r.Style = "Normal"
$r = $ws.Range('E27')
$r.NumberFormat = "@"
$r.Value = '  -5.11%  '
$r.Style = "Normal"
$r = $ws.Range('D28')
$r.NumberFormat = "@"
$r.Value = '1.00'
$r.Style = "Normal"
$r = $ws.Range('E28')
$r.NumberFormat = "@"
$r.Value = '  +0.05%  '
$r.Style = "Normal"
$r = $ws.Range('E29')
$r.NumberFormat = "@"
$r.Value = '  -2.95%  '
$r.Style = "Normal"
$r = $ws.Range('D30')
$r.NumberFormat = "@"
$r.Value = '22.87'
$r.Style = "Normal"
$r = $ws.Range('E30')
$r.NumberFormat = "@"
$r.Value = '  -2.73%  '
$r.Style = "Normal"
$r = $ws.Range('B31')
$r.NumberFormat = "@"
$r.Value = 'NEARProtocol'
$r.Style = "Normal"
$r = $ws.Range('C31')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$r.Style = "Normal"
$r = $ws.Range('D31')
$r.NumberFormat = "@"
$r.Value = '5.29'
$r.Style = "Normal"
$r = $ws.Range('E31')
$r.NumberFormat = "@"
$r.Value = '  -6.42%  '
$r.Style = "Normal"
$r = $ws.Range('B32')
$r.NumberFormat = "@"
$r.Value = 'USDe'
$r.Style = "Normal"
$r = $ws.Range('C32')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$r.Style = "Normal"
$r = $ws.Range('D32')
$r.NumberFormat = "@"
$r.Value = '0.999'
$r.Style = "Normal"
$r = $ws.Range('E32')
$r.NumberFormat = "@"
$r.Value = '  +0.01%  '
$r.Style = "Normal"
$r = $ws.Range('B33')
$r.NumberFormat = "@"
$r.Value = 'Fetch.AI'
$r.Style = "Normal"
$r = $ws.Range('C33')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$r.Style = "Normal"
$r = $ws.Range('D33')
$r.NumberFormat = "@"
$r.Value = '1.23'
$r.Style = "Normal"
$r = $ws.Range('E33')
$r.NumberFormat = "@"
$r.Value = '  -4.22%  '
$r.Style = "Normal"
$r = $ws.Range('B34')
$r.NumberFormat = "@"
$r.Value = 'Aptos'
$r.Style = "Normal"
$r = $ws.Range('C34')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$r.Style = "Normal"
$r = $ws.Range('D34')
$r.NumberFormat = "@"
$r.Value = '6.77'
$r.Style = "Normal"
$r = $ws.Range('E34')
$r.NumberFormat = "@"
$r.Value = '  -3.28%  '
$r.Style = "Normal"
$r = $ws.Range('D35')
$r.NumberFormat = "@"
$r.Value = '162.93'
$r.Style = "Normal"
$r = $ws.Range('E35')
$r.NumberFormat = "@"
$r.Value = '  -1.58%  '
$r.Style = "Normal"
$r = $ws.Range('E36')
$r.NumberFormat = "@"
$r.Value = '  -5.58%  '
$r.Style = "Normal"
$r = $ws.Range('D37')
$r.NumberFormat = "@"
$r.Value = '27.45'
$r.Style = "Normal"
$r = $ws.Range('E37')
$r.NumberFormat = "@"
$r.Value = '  -3.17%  '
$r.Style = "Normal"
$r = $ws.Range('D38')
$r.NumberFormat = "@"
$r.Value = '1.82'
$r.Style = "Normal"
$r = $ws.Range('E38')
$r.NumberFormat = "@"
$r.Value = '  -6.30%  '
$r.Style = "Normal"
$r = $ws.Range('D39')
$r.NumberFormat = "@"
$r.Value = '2.822.41'
$r.Style = "Normal"
$r = $ws.Range('E39')
$r.NumberFormat = "@"
$r.Value = '  +2.35%  '
$r.Style = "Normal"
$r = $ws.Range('D40')
$r.NumberFormat = "@"
$r.Value = '0.794'
$r.Style = "Normal"
$r = $ws.Range('E40')
$r.NumberFormat = "@"
$r.Value = '  -2.68%  '
$r.Style = "Normal"
$r = $ws.Range('E41')
$r.NumberFormat = "@"
$r.Value = '  -3.94%  '
$r.Style = "Normal"
$r = $ws.Range('E42')
$r.NumberFormat = "@"
$r.Value = '  -6.29%  '
$r.Style = "Normal"
$r = $ws.Range('D43')
$r.NumberFormat = "@"
$r.Value = '40.18'
$r.Style = "Normal"
$r = $ws.Range('E43')
$r.NumberFormat = "@"
$r.Value = '  -2.34%  '
$r.Style = "Normal"
$r = $ws.Range('E44')
$r.NumberFormat = "@"
$r.Value = '  -3.39%  '
$r.Style = "Normal"
$r = $ws.Range('D45')
$r.NumberFormat = "@"
$r.Value = '24.43'
$r.Style = "Normal"
$r = $ws.Range('E45')
$r.NumberFormat = "@"
$r.Value = '  -4.46%  '
$r.Style = "Normal"
$r = $ws.Range('E46')
$r.NumberFormat = "@"
$r.Value = '  -6.85%  '
$r.Style = "Normal"
$r = $ws.Range('D47')
$r.NumberFormat = "@"
$r.Value = '321.20'
$r.Style = "Normal"
$r = $ws.Range('E47')
$r.NumberFormat = "@"
$r.Value = '  -5.75%  '
$r.Style = "Normal"
$r = $ws.Range('D48')
$r.NumberFormat = "@"
$r.Value = '0.0274'
$r.Style = "Normal"
$r = $ws.Range('E48')
$r.NumberFormat = "@"
$r.Value = '  -4.03%  '
$r.Style = "Normal"
$r = $ws.Range('E49')
$r.NumberFormat = "@"
$r.Value = '  +1.12%  '
$r.Style = "Normal"
$r = $ws.Range('D50')
$r.NumberFormat = "@"
$r.Value = '0.978'
$r.Style = "Normal"
$r = $ws.Range('E50')
$r.NumberFormat = "@"
$r.Value = '  -4.24%  '
$r.Style = "Normal"
$r = $ws.Range('D51')
$r.NumberFormat = "@"
$r.Value = '6.16'
$r.Style = "Normal"
$r = $ws.Range('E51')
$r.NumberFormat = "@"
$r.Value = '  -2.68%  '
$r.Style = "Normal"
